$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.Value = "'" + $text
    $cellRange.ClearFormats()
}

Set-TextValue $ws.Range("D2") "44.019.13"
Set-TextValue $ws.Range("E2") "  -0.99%  "
Set-TextValue $ws.Range("D3") "2.240.82"
Set-TextValue $ws.Range("E3") "  -2.02%  "
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "314.99"
Set-TextValue $ws.Range("E5") "  -1.79%  "
Set-TextValue $ws.Range("D6") "99.10"
Set-TextValue $ws.Range("E6") "  -8.24%  "
Set-TextValue $ws.Range("D7") "0.574"
Set-TextValue $ws.Range("E7") "  -3.32%  "
Set-TextValue $ws.Range("E8") "  +0.17%  "
Set-TextValue $ws.Range("D9") "0.533"
Set-TextValue $ws.Range("E9") "  -7.72%  "
Set-TextValue $ws.Range("D10") "36.20"
Set-TextValue $ws.Range("E10") "  -7.57%  "
Set-TextValue $ws.Range("E11") "  -2.79%  "
Set-TextValue $ws.Range("D12") "7.37"
Set-TextValue $ws.Range("E12") "  -7.67%  "
Set-TextValue $ws.Range("E13") "  -3.06%  "
Set-TextValue $ws.Range("D14") "2.581.64"
Set-TextValue $ws.Range("E14") "  -2.04%  "
Set-TextValue $ws.Range("D15") "0.843"
Set-TextValue $ws.Range("E15") "  -5.52%  "
Set-TextValue $ws.Range("D16") "2.247.43"
Set-TextValue $ws.Range("E16") "  -1.92%  "
Set-TextValue $ws.Range("D17") "13.94"
Set-TextValue $ws.Range("E17") "  -5.47%  "
Set-TextValue $ws.Range("D18") "43.840.63"
Set-TextValue $ws.Range("E18") "  -1.18%  "
Set-TextValue $ws.Range("D19") "13.24"
Set-TextValue $ws.Range("E19") "  -6.92%  "
Set-TextValue $ws.Range("D20") "0.0₃0974"
Set-TextValue $ws.Range("E20") "  -3.22%  "
Set-TextValue $ws.Range("D21") "6.32"
Set-TextValue $ws.Range("D22") "65.99"
Set-TextValue $ws.Range("E22") "  -1.21%  "
Set-TextValue $ws.Range("D23") "237.79"
Set-TextValue $ws.Range("E23") "  -0.98%  "
Set-TextValue $ws.Range("E24") "  -8.18%  "
Set-TextValue $ws.Range("E25") "  -9.65%  "
Set-TextValue $ws.Range("E26") "  +0.16%  "
Set-TextValue $ws.Range("D27") "10.14"
Set-TextValue $ws.Range("E27") "  -1.77%  "
Set-TextValue $ws.Range("E28") "  -4.32%  "
Set-TextValue $ws.Range("D29") "36.51"
Set-TextValue $ws.Range("E29") "  -6.84%  "
Set-TextValue $ws.Range("D30") "5.99"
Set-TextValue $ws.Range("E30") "  -9.28%  "
Set-TextValue $ws.Range("E31") "  -3.87%  "
Set-TextValue $ws.Range("D32") "154.45"
Set-TextValue $ws.Range("E32") "  -5.63%  "
Set-TextValue $ws.Range("E33") "  -6.62%  "
Set-TextValue $ws.Range("E34") "  +1.08%  "
Set-TextValue $ws.Range("D35") "2.66"
Set-TextValue $ws.Range("E35") "  -2.53%  "
Set-TextValue $ws.Range("D36") "1.91"
Set-TextValue $ws.Range("E36") "  -8.19%  "
Set-TextValue $ws.Range("E37") "  -7.26%  "
Set-TextValue $ws.Range("D38") "0.118"
Set-TextValue $ws.Range("E38") "  -3.43%  "
Set-TextValue $ws.Range("D39") "15.75"
Set-TextValue $ws.Range("E39") "  +1.22%  "
Set-TextValue $ws.Range("D42") "0.0308"
Set-TextValue $ws.Range("E42") "  -6.99%  "
Set-TextValue $ws.Range("E43") "  +0.07%  "
Set-TextValue $ws.Range("D44") "1.705.84"
Set-TextValue $ws.Range("E44") "  -4.12%  "
Set-TextValue $ws.Range("D45") "82.25"
Set-TextValue $ws.Range("E45") "  -4.91%  "
Set-TextValue $ws.Range("E46") "  -6.96%  "
Set-TextValue $ws.Range("D47") "5.17"
Set-TextValue $ws.Range("E47") "  -5.98%  "
Set-TextValue $ws.Range("D48") "101.74"
Set-TextValue $ws.Range("E48") "  -3.43%  "
Set-TextValue $ws.Range("D49") "71.68"
Set-TextValue $ws.Range("E49") "  -5.69%  "
Set-TextValue $ws.Range("D50") "56.32"
Set-TextValue $ws.Range("E50") "  -7.18%  "
Set-TextValue $ws.Range("E51") "  -6.13%  "

# Rows 40 and 41 swap (RenderToken <-> NEARProtocol) with updated values
Set-TextValue $ws.Range("B40") "NEARProtocol"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D40") "3.52"
Set-TextValue $ws.Range("E40") "  -12.50%  "
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D41") "4.01"
Set-TextValue $ws.Range("E41") "  -11.22%  "
